$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-31 01:26:16"

$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
